# Generate Report for Handback
# Update the "Latest Handback DateTime" (column K) for the first data row
# (the 7990fd38-cd2c-40ec-a0bf-423934b0b2cd.md file) on both locale sheets,
# reflecting a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-11-15 17:45:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-15 17:46:13"
